# Update the last_edited_time values (column D) for the Notion "Chấm công HỆ THỐNG" export.
# Originally rows 2-14 shared the "2024-08-03T03:17:00.000Z" timestamp text, while rows
# 15-22 shared "2024-08-03T03:18:00.000Z". This edit bumps both timestamps by 11 minutes
# and also moves row 15 into the earlier ("...T03:28:00.000Z") bucket so its
# last_edited_time now matches rows 2-14 instead of rows 16-22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 2-14: "2024-08-03T03:17:00.000Z" -> "2024-08-03T03:28:00.000Z"
for ($r = 2; $r -le 14; $r++) {
    $ws.Range("D$r").Value = "2024-08-03T03:28:00.000Z"
}

# Row 15 moves from the "...T03:18:00.000Z" bucket to the "...T03:28:00.000Z" bucket
$ws.Range("D15").Value = "2024-08-03T03:28:00.000Z"

# Rows 16-22: "2024-08-03T03:18:00.000Z" -> "2024-08-03T03:29:00.000Z"
for ($r = 16; $r -le 22; $r++) {
    $ws.Range("D$r").Value = "2024-08-03T03:29:00.000Z"
}
